$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old extra columns (E:K) that are no longer used
$ws.Range("E1:K2").ClearContents()

# Update header row (A1:D1)
$ws.Range("A1").Value = "Day-Time"
$ws.Range("B1").Value = "Hall 1"
$ws.Range("C1").Value = "Hall 2"
$ws.Range("D1").Value = "Hall 5"

# Update data row (A2:D2)
$ws.Range("A2").Value = "9AM - 12PM"
$ws.Range("B2").Value = "CD222"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"

$wb.Save()
